$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.508.07'
$ws.Range("E2").Value = '  +2.36%  '

# Row 3
$ws.Range("D3").Value = '3.393.02'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.33'
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.75%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.597'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.95%  '

# Row 9
$ws.Range("E9").Value = '  +5.77%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.593'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.56%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.63'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.54%  '

# Row 12
$ws.Range("E12").Value = '  +3.06%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '681.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.67%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.16%  '

# Row 15
$ws.Range("D15").Value = '3.935.71'
$ws.Range("E15").Value = '  +1.70%  '

# Row 16
$ws.Range("D16").Value = '69.504.19'
$ws.Range("E16").Value = '  +2.26%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.395.15'
$ws.Range("E17").Value = '  +1.85%  '

# Row 18
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.120'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.71%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.95%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.78%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.905'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.76%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.54%  '

# Row 23
$ws.Range("E23").Value = '  +0.65%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.54%  '

# Row 25
$ws.Range("E25").Value = '  +0.56%  '

# Row 26
$ws.Range("E26").Value = '  +1.09%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.17%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '34.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.16%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.44%  '

# Row 30
$ws.Range("E30").Value = '  -0.45%  '

# Row 31
$ws.Range("E31").Value = '  +0.91%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '559.55'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.16%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.64%  '

# Row 34
$ws.Range("E34").Value = '  +1.02%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.84%  '

# Row 36
$ws.Range("E36").Value = '  +0.03%  '

# Row 37
$ws.Range("D37").Value = '3.693.75'
$ws.Range("E37").Value = '  +0.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.60'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.37%  '

# Row 39
$ws.Range("E39").Value = '  +4.67%  '

# Row 40
$ws.Range("E40").Value = '  +3.40%  '

# Row 41
$ws.Range("E41").Value = '  +1.86%  '

# Row 42
$ws.Range("E42").Value = '  +3.03%  '

# Row 43
$ws.Range("E43").Value = '  +0.62%  '

# Row 44
$ws.Range("E44").Value = '  +4.12%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.79%  '

# Row 46
$ws.Range("E46").Value = '  -0.17%  '

# Row 47
$ws.Range("E47").Value = '  +0.75%  '

# Row 48
$ws.Range("E48").Value = '  +5.05%  '

# Row 49
$ws.Range("E49").Value = '  -0.10%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.05'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.48%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.38%  '
